$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("M1").Value = "renewd"
$ws.Range("N1").Value = "PlanID"
$ws.Range("O1").Value = "iteration"

# Copy header style (bold/centered/bordered) from an existing header cell to the new ones
$ws.Range("L1").Copy() | Out-Null
$ws.Range("M1:O1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill data rows 2..31
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 13).Value = "after"       # column M
    $ws.Cells.Item($r, 14).Value = 20140060      # column N
    $ws.Cells.Item($r, 15).Value = 11            # column O
}
